$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

$cell = $tbl.Cell(1, 1)
$cell.Range.Text = "40+56=96"
$cell = $tbl.Cell(1, 2)
$cell.Range.Text = "9+76=85"
$cell = $tbl.Cell(1, 3)
$cell.Range.Text = "14+84=98"
$cell = $tbl.Cell(1, 4)
$cell.Range.Text = "52-42=10"
$cell = $tbl.Cell(1, 5)
$cell.Range.Text = "85-54=31"
$cell = $tbl.Cell(2, 1)
$cell.Range.Text = "54-23=31"
$cell = $tbl.Cell(2, 2)
$cell.Range.Text = "0+63=63"
$cell = $tbl.Cell(2, 3)
$cell.Range.Text = "12+16=28"
$cell = $tbl.Cell(2, 4)
$cell.Range.Text = "39-34=5"
$cell = $tbl.Cell(2, 5)
$cell.Range.Text = "95-1=94"
$cell = $tbl.Cell(3, 1)
$cell.Range.Text = "80-8=72"
$cell = $tbl.Cell(3, 2)
$cell.Range.Text = "20+60=80"
$cell = $tbl.Cell(3, 3)
$cell.Range.Text = "24-18=6"
$cell = $tbl.Cell(3, 4)
$cell.Range.Text = "83-21=62"
$cell = $tbl.Cell(3, 5)
$cell.Range.Text = "11+74=85"
$cell = $tbl.Cell(4, 1)
$cell.Range.Text = "81-25=56"
$cell = $tbl.Cell(4, 2)
$cell.Range.Text = "5+3=8"
$cell = $tbl.Cell(4, 3)
$cell.Range.Text = "33+66=99"
$cell = $tbl.Cell(4, 4)
$cell.Range.Text = "92-32=60"
$cell = $tbl.Cell(4, 5)
$cell.Range.Text = "20-0=20"
$cell = $tbl.Cell(5, 1)
$cell.Range.Text = "36-6=30"
$cell = $tbl.Cell(5, 2)
$cell.Range.Text = "80-0=80"
$cell = $tbl.Cell(5, 3)
$cell.Range.Text = "7-5=2"
$cell = $tbl.Cell(5, 4)
$cell.Range.Text = "62-39=23"
$cell = $tbl.Cell(5, 5)
$cell.Range.Text = "84-30=54"
$cell = $tbl.Cell(6, 1)
$cell.Range.Text = "7+32=39"
$cell = $tbl.Cell(6, 2)
$cell.Range.Text = "7-3=4"
$cell = $tbl.Cell(6, 3)
$cell.Range.Text = "71-33=38"
$cell = $tbl.Cell(6, 4)
$cell.Range.Text = "40-8=32"
$cell = $tbl.Cell(6, 5)
$cell.Range.Text = "59-15=44"
$cell = $tbl.Cell(7, 1)
$cell.Range.Text = "71-12=59"
$cell = $tbl.Cell(7, 2)
$cell.Range.Text = "7+89=96"
$cell = $tbl.Cell(7, 3)
$cell.Range.Text = "97-32=65"
$cell = $tbl.Cell(7, 4)
$cell.Range.Text = "12+1=13"
$cell = $tbl.Cell(7, 5)
$cell.Range.Text = "73+4=77"
$cell = $tbl.Cell(8, 1)
$cell.Range.Text = "18+3=21"
$cell = $tbl.Cell(8, 2)
$cell.Range.Text = "50+36=86"
$cell = $tbl.Cell(8, 3)
$cell.Range.Text = "53+32=85"
$cell = $tbl.Cell(8, 4)
$cell.Range.Text = "35-33=2"
$cell = $tbl.Cell(8, 5)
$cell.Range.Text = "60+13=73"
$cell = $tbl.Cell(9, 1)
$cell.Range.Text = "5+71=76"
$cell = $tbl.Cell(9, 2)
$cell.Range.Text = "91-8=83"
$cell = $tbl.Cell(9, 3)
$cell.Range.Text = "55+11=66"
$cell = $tbl.Cell(9, 4)
$cell.Range.Text = "11+55=66"
$cell = $tbl.Cell(9, 5)
$cell.Range.Text = "12-1=11"
$cell = $tbl.Cell(10, 1)
$cell.Range.Text = "16+45=61"
$cell = $tbl.Cell(10, 2)
$cell.Range.Text = "35-15=20"
$cell = $tbl.Cell(10, 3)
$cell.Range.Text = "20+64=84"
$cell = $tbl.Cell(10, 4)
$cell.Range.Text = "2+51=53"
$cell = $tbl.Cell(10, 5)
$cell.Range.Text = "35+25=60"
$cell = $tbl.Cell(11, 1)
$cell.Range.Text = "92+4=96"
$cell = $tbl.Cell(11, 2)
$cell.Range.Text = "90-70=20"
$cell = $tbl.Cell(11, 3)
$cell.Range.Text = "9+31=40"
$cell = $tbl.Cell(11, 4)
$cell.Range.Text = "73-57=16"
$cell = $tbl.Cell(11, 5)
$cell.Range.Text = "43-3=40"
$cell = $tbl.Cell(12, 1)
$cell.Range.Text = "33+32=65"
$cell = $tbl.Cell(12, 2)
$cell.Range.Text = "45-0=45"
$cell = $tbl.Cell(12, 3)
$cell.Range.Text = "31+42=73"
$cell = $tbl.Cell(12, 4)
$cell.Range.Text = "20+72=92"
$cell = $tbl.Cell(12, 5)
$cell.Range.Text = "83-28=55"
$cell = $tbl.Cell(13, 1)
$cell.Range.Text = "63-35=28"
$cell = $tbl.Cell(13, 2)
$cell.Range.Text = "51-21=30"
$cell = $tbl.Cell(13, 3)
$cell.Range.Text = "81-0=81"
$cell = $tbl.Cell(13, 4)
$cell.Range.Text = "22+59=81"
$cell = $tbl.Cell(13, 5)
$cell.Range.Text = "49-45=4"
$cell = $tbl.Cell(14, 1)
$cell.Range.Text = "90-51=39"
$cell = $tbl.Cell(14, 2)
$cell.Range.Text = "46-32=14"
$cell = $tbl.Cell(14, 3)
$cell.Range.Text = "40-34=6"
$cell = $tbl.Cell(14, 4)
$cell.Range.Text = "83-26=57"
$cell = $tbl.Cell(14, 5)
$cell.Range.Text = "42+8=50"
$cell = $tbl.Cell(15, 1)
$cell.Range.Text = "54-30=24"
$cell = $tbl.Cell(15, 2)
$cell.Range.Text = "83-42=41"
$cell = $tbl.Cell(15, 3)
$cell.Range.Text = "24+35=59"
$cell = $tbl.Cell(15, 4)
$cell.Range.Text = "68-67=1"
$cell = $tbl.Cell(15, 5)
$cell.Range.Text = "86-38=48"
$cell = $tbl.Cell(16, 1)
$cell.Range.Text = "24+58=82"
$cell = $tbl.Cell(16, 2)
$cell.Range.Text = "10+58=68"
$cell = $tbl.Cell(16, 3)
$cell.Range.Text = "3+67=70"
$cell = $tbl.Cell(16, 4)
$cell.Range.Text = "17+1=18"
$cell = $tbl.Cell(16, 5)
$cell.Range.Text = "92-67=25"
$cell = $tbl.Cell(17, 1)
$cell.Range.Text = "66-39=27"
$cell = $tbl.Cell(17, 2)
$cell.Range.Text = "1+28=29"
$cell = $tbl.Cell(17, 3)
$cell.Range.Text = "15+78=93"
$cell = $tbl.Cell(17, 4)
$cell.Range.Text = "24+59=83"
$cell = $tbl.Cell(17, 5)
$cell.Range.Text = "17+69=86"
$cell = $tbl.Cell(18, 1)
$cell.Range.Text = "59-6=53"
$cell = $tbl.Cell(18, 2)
$cell.Range.Text = "80-74=6"
$cell = $tbl.Cell(18, 3)
$cell.Range.Text = "30-19=11"
$cell = $tbl.Cell(18, 4)
$cell.Range.Text = "77-9=68"
$cell = $tbl.Cell(18, 5)
$cell.Range.Text = "72-34=38"
$cell = $tbl.Cell(19, 1)
$cell.Range.Text = "75-32=43"
$cell = $tbl.Cell(19, 2)
$cell.Range.Text = "4+31=35"
$cell = $tbl.Cell(19, 3)
$cell.Range.Text = "23+2=25"
$cell = $tbl.Cell(19, 4)
$cell.Range.Text = "68+2=70"
$cell = $tbl.Cell(19, 5)
$cell.Range.Text = "6+28=34"
$cell = $tbl.Cell(20, 1)
$cell.Range.Text = "60+35=95"
$cell = $tbl.Cell(20, 2)
$cell.Range.Text = "15+35=50"
$cell = $tbl.Cell(20, 3)
$cell.Range.Text = "79-30=49"
$cell = $tbl.Cell(20, 4)
$cell.Range.Text = "28+68=96"
$cell = $tbl.Cell(20, 5)
$cell.Range.Text = "66-38=28"

Write-Output "done"
